# Fruta / hortaliza, semanal
# Insert a new weekly record at row 3 (pushing the existing data rows down by one),
# then populate the new row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; this shifts former rows 3..29 down to 4..30
# and inherits formatting (including the date style on column D) from the row above.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Cells.Item(3, 1).Value2  = 11
$ws.Cells.Item(3, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value2  = "Bíobío"
$ws.Cells.Item(3, 4).Value2  = 44530
$ws.Cells.Item(3, 5).Value2  = 8
$ws.Cells.Item(3, 6).Value2  = 100112031
$ws.Cells.Item(3, 7).Value2  = "Poroto verde"
$ws.Cells.Item(3, 8).Value2  = "Sin especificar"
$ws.Cells.Item(3, 9).Value2  = "Primera"
$ws.Cells.Item(3, 10).Value2 = 110
$ws.Cells.Item(3, 11).Value2 = 22000
$ws.Cells.Item(3, 12).Value2 = 23000
$ws.Cells.Item(3, 13).Value2 = 22455
$ws.Cells.Item(3, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(3, 15).Value2 = "Región del Maule"
$ws.Cells.Item(3, 16).Value2 = 898
$ws.Cells.Item(3, 17).Value2 = 25
$ws.Cells.Item(3, 18).Value2 = "Hortaliza"
